$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.463
$ws.Range("A8").Value = -21.255
$ws.Range("A10").Value = -21.047
$ws.Range("A12").Value = -21.808
$ws.Range("B13").Value = 6.724000000000001
$ws.Range("A18").Value = -21.808
$ws.Range("C20").Value = -13.041
